$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header format (G1, e.g. "sum") onto the new H1 header cell
# so the new "Save" column header matches the look of the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the new "Save" column header
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for the data rows (2-5)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
